$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; this shifts rows 15-77 down to 16-78
$ws.Rows.Item(15).EntireRow.Insert()

# Populate the newly inserted row 15 with its values
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value = "Maule"
$ws.Cells.Item(15, 4).Value = 44592
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103002
$ws.Cells.Item(15, 10).Value = "Ciruela"
$ws.Cells.Item(15, 11).Value = "Black Amber"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 400
$ws.Cells.Item(15, 14).Value = 9000
$ws.Cells.Item(15, 15).Value = 9000
$ws.Cells.Item(15, 16).Value = 9000
$ws.Cells.Item(15, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(15, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(15, 19).Value = 500
$ws.Cells.Item(15, 20).Value = 18
